# Durham VOC positivity workbook update:
# A new day of data (date 44285, serial for 2021-04-17) was added at the top
# of the data table on "Sheet2". Insert a new row 2 (pushing the existing
# data and the footer rows down by one) and populate it with the new day's
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2; existing rows 2-39 shift down to 3-40,
# formulas / relative refs / styles shift automatically with them.
$ws.Rows("2:2").Insert()

# Fill in the new row with the latest day's values.
$ws.Range("A2").Value = 44285
$ws.Range("A2").NumberFormat = "d-mmm"

$ws.Range("B2").Value = 0.1398
$ws.Range("C2").Value = 0.1765
$ws.Range("D2").Value = 0.1091
$ws.Range("E2").Value = 0.1207
$ws.Range("F2").Value = 0.1102
$ws.Range("G2").Value = 0.061
$ws.Range("H2").Value = 0.0686
$ws.Range("I2").Value = 0.0941

$ws.Range("J2").Value = 44285
$ws.Range("K2").Value = 31

# Match the cursor position left behind by the author after editing.
$ws.Range("P28").Select()
